# Update the "Periodo Mora" column (E16:E43) so the periods run in
# ascending order (1607 .. 2003) instead of the previous descending
# order (2003 .. 1607). Row styling / borders stay exactly as they were;
# only the text values change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$periods = @(
    "1607",
    "1608",
    "1802",
    "1803",
    "1804",
    "1805",
    "1806",
    "1807",
    "1808",
    "1809",
    "1810",
    "1811",
    "1812",
    "1901",
    "1902",
    "1903",
    "1904",
    "1905",
    "1906",
    "1907",
    "1908",
    "1909",
    "1910",
    "1911",
    "1912",
    "2001",
    "2002",
    "2003"
)

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
}
